# Add a new column F for "pvt1" (solar_th / pvt), mirroring the layout
# already used for the other asset columns (net1, pv1, bat1, CHP1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in F1, formatted like the other header cells (B1:E1).
$ws.Range("F1").Value = "pvt1"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

# New data cells F2:F4, matching the plain (unstyled) numeric cells below
# the header row.
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0

$excel.CutCopyMode = 0
